# Income_details workbook edit.
#
# Before:
#   Row2: A2="bonus"   B2=2000  C2=46070.250231481485
#   Row3: A3="salary2" B3=1000  C3=46069.250231481485
#
# After:
#   Row2: (A2 cleared)  B2=3000   C2=46071.250231481485
#   Row3: (A3 cleared)  B3=7000   C3=46070.250231481485
#   Row4: (new)         B4=10000  C4=46064.250231481485  (same date style as C2/C3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the text labels that used to live in column A for rows 2 and 3 -
# the refreshed rows no longer carry a "Source" label.
$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()

# Update the amount/date pair for row 2.
$ws.Range("B2").Value = 3000
$ws.Range("C2").Value = 46071.250231481485

# Update the amount/date pair for row 3.
$ws.Range("B3").Value = 7000
$ws.Range("C3").Value = 46070.250231481485

# Add a brand new row 4 with its own amount/date pair. Copy the date cell's
# formatting from C2 first so C4 picks up the same date style (rather than
# minting a brand-new number format/style entry), then overwrite the values.
$ws.Range("C2").Copy($ws.Range("C4"))
$ws.Range("B4").Value = 10000
$ws.Range("C4").Value = 46064.250231481485
